$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6122626666666666
$ws.Range("H2").Value = 1.836788
$ws.Range("I2").Value = 0.006779070576782467
$ws.Range("J2").Value = 0.006779070576782467
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 0.03180173925688889
$ws.Range("R2").Value = 0.2862156533119999
$ws.Range("S2").Value = 0.0001879848396509225
$ws.Range("T2").Value = 0.0001879848396509225
$ws.Range("G3").Value = 0.6122626666666666
$ws.Range("H3").Value = 1.836788
$ws.Range("I3").Value = 0.006779070576782467
$ws.Range("J3").Value = 0.006779070576782467
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 1.115026033063556
$ws.Range("R3").Value = 10.035234297572
$ws.Range("S3").Value = 0.006591085737131544
$ws.Range("T3").Value = 0.006591085737131545
$ws.Range("I4").Value = 0.003538518590750013
$ws.Range("J4").Value = 0.003538518590750013
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("S4").Value = 0.00009812375344817044
$ws.Range("T4").Value = 0.00009812375344817044
$ws.Range("I5").Value = 0.003538518590750013
$ws.Range("J5").Value = 0.003538518590750013
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 0.5820178891010001
$ws.Range("R5").Value = 5.238161001909
$ws.Range("S5").Value = 0.003440394837301843
$ws.Range("T5").Value = 0.003440394837301843
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.115957
$ws.Range("H6").Value = 0.347871
$ws.Range("I6").Value = 0.001283894527085267
$ws.Range("J6").Value = 0.001283894527085267
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.006022961189333333
$ws.Range("R6").Value = 0.05420665070399999
$ws.Range("S6").Value = 0.00003560262488333224
$ws.Range("T6").Value = 0.00003560262488333224
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.115957
$ws.Range("H7").Value = 0.347871
$ws.Range("I7").Value = 0.001283894527085267
$ws.Range("J7").Value = 0.001283894527085267
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 0.2111758249443333
$ws.Range("R7").Value = 1.900582424499
$ws.Range("S7").Value = 0.001248291902201935
$ws.Range("T7").Value = 0.001248291902201935
$ws.Range("G8").Value = 89.26880233333334
$ws.Range("H8").Value = 267.806407
$ws.Range("I8").Value = 0.9883985163053822
$ws.Range("J8").Value = 0.9883985163053823
$ws.Range("O8").Value = 0.02773017886769741
$ws.Range("P8").Value = 0.02773017886769741
$ws.Range("Q8").Value = 4.636740618263111
$ws.Range("R8").Value = 41.730665564368
$ws.Range("S8").Value = 0.02740846764971499
$ws.Range("T8").Value = 0.02740846764971499
$ws.Range("G9").Value = 89.26880233333334
$ws.Range("H9").Value = 267.806407
$ws.Range("I9").Value = 0.9883985163053822
$ws.Range("J9").Value = 0.9883985163053823
$ws.Range("M9").Value = 1.821156333333333
$ws.Range("N9").Value = 5.463469
$ws.Range("O9").Value = 0.9722698211323025
$ws.Range("P9").Value = 0.9722698211323026
$ws.Range("Q9").Value = 162.5724447384315
$ws.Range("R9").Value = 1463.152002645883
$ws.Range("S9").Value = 0.9609900486556672
$ws.Range("T9").Value = 0.9609900486556674
